$wb = $excel.ActiveWorkbook

# --- newApp_ApplicationDetailsMuraba: update record_reference_number (column C, row 6) ---
$wsMuraba = $wb.Worksheets.Item("newApp_ApplicationDetailsMuraba")
$wsMuraba.Range("C6").Value = "4633"

# --- TW_NewApp_AppDetails_TestData: update record_reference_number (column C, row 6) ---
$wsTawarruq = $wb.Worksheets.Item("TW_NewApp_AppDetails_TestData")
$wsTawarruq.Range("C6").Value = "4635"

# --- ULSExecution: update Test Execution Status column (C) for rows 2-123 ---
$wsExec = $wb.Worksheets.Item("ULSExecution")
for ($r = 2; $r -le 123; $r++) {
    if ($r -eq 3) {
        $wsExec.Range("C3").Value = "Passed"
    } else {
        $wsExec.Cells.Item($r, 3).Value = "Failed"
    }
}
